# Auto-generated Excel COM-interop script
# Applies updated crypto price/volume data (coin list refresh) to sheet1
# Source: commit 'Updated cryptos list on Fri May 12 20:42:37 UTC 2023 with GitHub Actions'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.475.47'
$ws.Range('E2').Value = '  -1.88%  '
$ws.Range('D3').Value = '1.788.32'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = "'308.80"
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = "'1.003"
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = "'0.4262"
$ws.Range('E7').Value = '  +1.77%  '
$ws.Range('D8').Value = "'0.3624"
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('D9').Value = "'0.07153"
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('D10').Value = "'0.8523"
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('D11').Value = "'20.53"
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('D12').Value = '1.824.45'
$ws.Range('E12').Value = '  +1.61%  '
$ws.Range('D13').Value = "'6.501"
$ws.Range('E13').Value = '  +2.46%  '
$ws.Range('D14').Value = "'5.263"
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = "'0.06881"
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = "'1.003"
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = "'79.71"
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = "'0.000008858"
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').Value = "'15.00"
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').Value = '26.502.34'
$ws.Range('E21').Value = '  -1.83%  '
$ws.Range('D22').Value = "'5.139"
$ws.Range('E22').Value = '  +1.76%  '
$ws.Range('D23').Value = "'11.02"
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('D24').Value = '2.054.47'
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('D25').Value = "'152.03"
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('D26').Value = "'1.821"
$ws.Range('E26').Value = '  -5.81%  '
$ws.Range('D27').Value = "'18.14"
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('D28').Value = "'5.146"
$ws.Range('E28').Value = '  +3.06%  '
$ws.Range('D29').Value = "'1.897"
$ws.Range('E29').Value = '  +15.82%  '
$ws.Range('D30').Value = "'114.66"
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('D31').Value = "'0.08888"
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('D32').Value = "'0.7448"
$ws.Range('E32').Value = '  +3.87%  '
$ws.Range('D33').Value = "'1.142"
$ws.Range('E33').Value = '  +6.35%  '
$ws.Range('D34').Value = "'4.346"
$ws.Range('E34').Value = '  +1.45%  '
$ws.Range('D35').Value = "'2.761"
$ws.Range('E35').Value = '  -3.49%  '
$ws.Range('D36').Value = "'1.003"
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = "'1.117"
$ws.Range('E37').Value = '  +4.15%  '
$ws.Range('D38').Value = "'0.05152"
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('D39').Value = "'0.01900"
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = "'0.4981"
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('D41').Value = "'0.1620"
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').Value = "'2.602"
$ws.Range('E42').Value = '  +1.12%  '
$ws.Range('D43').Value = "'6.404"
$ws.Range('E43').Value = '  +7.60%  '
$ws.Range('D44').Value = "'8.237"
$ws.Range('E44').Value = '  +2.78%  '
$ws.Range('D45').Value = "'105.55"
$ws.Range('E45').Value = '  +1.17%  '
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('D47').Value = "'1.003"
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').Value = "'1.637"
$ws.Range('E48').Value = '  +2.58%  '
$ws.Range('D49').Value = "'0.4520"
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('D50').Value = "'0.06202"
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('D51').Value = "'1.761"
$ws.Range('E51').Value = '  +4.63%  '
